$wb = $excel.ActiveWorkbook

# --- zh-cn sheet: row 7 (aec8b0ae-c7da-4752-8c00-5711c82b056d) ---
$wsZh = $wb.Worksheets.Item("zh-cn")

$i7 = $wsZh.Range("I7")
$i7.Value = "aec8b0ae-c7da-4752-8c00-5711c82b056d.md"
$i7.Font.Color = 15570276
$i7.Font.Underline = 2

$wsZh.Range("J7").Value = "aec8b0ae-c7da-4752-8c00-5711c82b056d.39e2e2e7ef886083ca05f0618dbd76776ca16f08.zh-cn.xlf"
$wsZh.Range("K7").Value = "2016-08-21 12:59:47"
$wsZh.Range("P7").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ff6f9bb334632ceea5742b1d2f752ee36447aeb0/e2e/aec8b0ae-c7da-4752-8c00-5711c82b056d.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/af61029ca465c17da26c298d29d2f70b7556923d/e2e/aec8b0ae-c7da-4752-8c00-5711c82b056d.md."

$wsZh.Hyperlinks.Add($i7, "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/af61029ca465c17da26c298d29d2f70b7556923d/e2e/aec8b0ae-c7da-4752-8c00-5711c82b056d.md", [Type]::Missing, [Type]::Missing, "aec8b0ae-c7da-4752-8c00-5711c82b056d.md")

# --- de-de sheet: rows 2-6 shift their shared-string indices (new strings were
# inserted before them), and row 7 gets the same treatment as zh-cn row 7 ---
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("G2").Value = "4b4c079f-0cbc-42bd-a447-1df24f3b4675.70f9d7b44d89b0d232697761fbbf661a292b2811.de-de.xlf"
$wsDe.Range("J2").Value = "4b4c079f-0cbc-42bd-a447-1df24f3b4675.70f9d7b44d89b0d232697761fbbf661a292b2811.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-21 12:54:19"

$wsDe.Range("G3").Value = "64399e0e-788c-4ed1-8cbc-6cf05aef8959.fc101161ee32c60ae45ce7fc08a046b09b90fc9a.de-de.xlf"
$wsDe.Range("J3").Value = "64399e0e-788c-4ed1-8cbc-6cf05aef8959.fc101161ee32c60ae45ce7fc08a046b09b90fc9a.de-de.xlf"
$wsDe.Range("K3").Value = "2016-08-21 12:56:43"

$wsDe.Range("G4").Value = "64399e0e-788c-4ed1-8cbc-6cf05aef8959.fc101161ee32c60ae45ce7fc08a046b09b90fc9a.de-de.xlf"
$wsDe.Range("J4").Value = "64399e0e-788c-4ed1-8cbc-6cf05aef8959.fc101161ee32c60ae45ce7fc08a046b09b90fc9a.de-de.xlf"
$wsDe.Range("K4").Value = "2016-08-21 12:56:43"

$wsDe.Range("G5").Value = "f7d892d6-6bc5-416f-8580-8922818b8172.e9c720208498288e90cccbf751f460269a1e9380.de-de.xlf"
$wsDe.Range("J5").Value = "f7d892d6-6bc5-416f-8580-8922818b8172.e9c720208498288e90cccbf751f460269a1e9380.de-de.xlf"
$wsDe.Range("K5").Value = "2016-08-21 12:57:45"

$wsDe.Range("G6").Value = "f5cf148b-fb95-41d2-9182-15a4abdcef62.d514dae71453899cae3fbae038f45b6bafa9ff08.de-de.xlf"
$wsDe.Range("P6").Value = "Handback file name: sbrqdifs.bza is different with handoff file name: f5cf148b-fb95-41d2-9182-15a4abdcef62.d514dae71453899cae3fbae038f45b6bafa9ff08.de-de."

$wsDe.Range("G7").Value = "aec8b0ae-c7da-4752-8c00-5711c82b056d.39e2e2e7ef886083ca05f0618dbd76776ca16f08.de-de.xlf"

$i7de = $wsDe.Range("I7")
$i7de.Value = "aec8b0ae-c7da-4752-8c00-5711c82b056d.md"
$i7de.Font.Color = 15570276
$i7de.Font.Underline = 2

$wsDe.Range("J7").Value = "aec8b0ae-c7da-4752-8c00-5711c82b056d.39e2e2e7ef886083ca05f0618dbd76776ca16f08.de-de.xlf"
$wsDe.Range("K7").Value = "2016-08-21 12:59:53"
$wsDe.Range("P7").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ff6f9bb334632ceea5742b1d2f752ee36447aeb0/e2e/aec8b0ae-c7da-4752-8c00-5711c82b056d.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/af61029ca465c17da26c298d29d2f70b7556923d/e2e/aec8b0ae-c7da-4752-8c00-5711c82b056d.md."

$wsDe.Hyperlinks.Add($i7de, "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/af61029ca465c17da26c298d29d2f70b7556923d/e2e/aec8b0ae-c7da-4752-8c00-5711c82b056d.md", [Type]::Missing, [Type]::Missing, "aec8b0ae-c7da-4752-8c00-5711c82b056d.md")
